$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Cells.Item(33, 8).Value = 96.25
$ws.Cells.Item(33, 9).Value = 108.70588
$ws.Cells.Item(33, 11).Value = 108.70588
$ws.Cells.Item(33, 13).Value = 120.29412

# Row 51
$ws.Cells.Item(51, 8).Value = 8000
$ws.Cells.Item(51, 9).Value = 8000
$ws.Cells.Item(51, 10).Value = 0
$ws.Cells.Item(51, 11).Value = 8000
$ws.Cells.Item(51, 12).Value = 0
$ws.Cells.Item(51, 13).Value = -7516
$ws.Cells.Item(51, 14).ClearContents()

# Row 99
$ws.Cells.Item(99, 8).Value = 631.25
$ws.Cells.Item(99, 9).Value = 610.4286
$ws.Cells.Item(99, 11).Value = 1831.2858
$ws.Cells.Item(99, 13).Value = -333.2857999999999

# Row 132
$ws.Cells.Item(132, 8).Value = 779.75
$ws.Cells.Item(132, 9).Value = 779.75
$ws.Cells.Item(132, 11).Value = 2339.25
$ws.Cells.Item(132, 13).Value = 190.75


$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 5645.2144
$ws.Cells.Item(32, 9).Value = 1507.591
$ws.Cells.Item(32, 10).Value = 20816.5
$ws.Cells.Item(32, 11).Value = 1507.591
$ws.Cells.Item(32, 12).Value = 20816.5
$ws.Cells.Item(32, 13).Value = -1220.591
$ws.Cells.Item(32, 14).Value = -21390.5

# Row 45
$ws.Cells.Item(45, 8).Value = 2571.5
$ws.Cells.Item(45, 9).Value = 2577
$ws.Cells.Item(45, 10).Value = 2566
$ws.Cells.Item(45, 11).Value = 2577
$ws.Cells.Item(45, 12).Value = 2566
$ws.Cells.Item(45, 13).Value = -2200
$ws.Cells.Item(45, 14).Value = -3320

# Row 61
$ws.Cells.Item(61, 8).Value = 2412.5
$ws.Cells.Item(61, 9).Value = 2344.4443
$ws.Cells.Item(61, 10).Value = 2616.6667
$ws.Cells.Item(61, 11).Value = 2344.4443
$ws.Cells.Item(61, 12).Value = 2616.6667
$ws.Cells.Item(61, 13).Value = -2132.4443
$ws.Cells.Item(61, 14).Value = -3040.6667

# Row 74
$ws.Cells.Item(74, 8).Value = 877.9545000000001
$ws.Cells.Item(74, 9).Value = 877.9545000000001
$ws.Cells.Item(74, 11).Value = 877.9545000000001
$ws.Cells.Item(74, 13).Value = -3.954500000000053

# Row 77
$ws.Cells.Item(77, 8).Value = 877.9545000000001
$ws.Cells.Item(77, 9).Value = 877.9545000000001
$ws.Cells.Item(77, 11).Value = 4389.7725
$ws.Cells.Item(77, 13).Value = -21.77250000000004

# Row 88
$ws.Cells.Item(88, 8).Value = 2997.4285
$ws.Cells.Item(88, 9).Value = 1387.5
$ws.Cells.Item(88, 10).Value = 3641.4
$ws.Cells.Item(88, 11).Value = 1387.5
$ws.Cells.Item(88, 12).Value = 3641.4
$ws.Cells.Item(88, 13).Value = -981.5
$ws.Cells.Item(88, 14).Value = -4453.4

# Row 91
$ws.Cells.Item(91, 8).Value = 2997.4285
$ws.Cells.Item(91, 9).Value = 1387.5
$ws.Cells.Item(91, 10).Value = 3641.4
$ws.Cells.Item(91, 11).Value = 1387.5
$ws.Cells.Item(91, 12).Value = 3641.4
$ws.Cells.Item(91, 13).Value = 16.5
$ws.Cells.Item(91, 14).Value = -6449.4

# Row 132
$ws.Cells.Item(132, 8).Value = 2806.1538
$ws.Cells.Item(132, 9).Value = 2545.6667
$ws.Cells.Item(132, 10).Value = 3900.2
$ws.Cells.Item(132, 11).Value = 7637.000100000001
$ws.Cells.Item(132, 12).Value = 11700.6
$ws.Cells.Item(132, 13).Value = -5107.000100000001
$ws.Cells.Item(132, 14).Value = -16760.6

# Row 136
$ws.Cells.Item(136, 8).Value = 2412.5
$ws.Cells.Item(136, 9).Value = 2344.4443
$ws.Cells.Item(136, 10).Value = 2616.6667
$ws.Cells.Item(136, 11).Value = 7033.3329
$ws.Cells.Item(136, 12).Value = 7850.000100000001
$ws.Cells.Item(136, 13).Value = -4483.3329
$ws.Cells.Item(136, 14).Value = -12950.0001


$ws = $wb.Worksheets.Item("BSM")
# Row 26
$ws.Cells.Item(26, 8).Value = 14333.333
$ws.Cells.Item(26, 9).Value = 14333.333
$ws.Cells.Item(26, 11).Value = 14333.333
$ws.Cells.Item(26, 13).Value = -14041.333

# Row 105
$ws.Cells.Item(105, 8).Value = 4929.6665
$ws.Cells.Item(105, 9).Value = 4929.6665
$ws.Cells.Item(105, 11).Value = 4929.6665
$ws.Cells.Item(105, 13).Value = -3182.6665

# Row 134
$ws.Cells.Item(134, 8).Value = 7890.5835
$ws.Cells.Item(134, 9).Value = 7299.3
$ws.Cells.Item(134, 11).Value = 21897.9
$ws.Cells.Item(134, 13).Value = -19362.9


$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 1178.7059
$ws.Cells.Item(31, 9).Value = 1092.5
$ws.Cells.Item(31, 10).Value = 1301.8572
$ws.Cells.Item(31, 11).Value = 1092.5
$ws.Cells.Item(31, 12).Value = 1301.8572
$ws.Cells.Item(31, 13).Value = -797.5
$ws.Cells.Item(31, 14).Value = -1891.8572

# Row 34
$ws.Cells.Item(34, 8).Value = 1178.7059
$ws.Cells.Item(34, 9).Value = 1092.5
$ws.Cells.Item(34, 10).Value = 1301.8572
$ws.Cells.Item(34, 11).Value = 1092.5
$ws.Cells.Item(34, 12).Value = 1301.8572
$ws.Cells.Item(34, 13).Value = -890.5
$ws.Cells.Item(34, 14).Value = -1705.8572

# Row 132
$ws.Cells.Item(132, 8).Value = 3882.1428
$ws.Cells.Item(132, 9).Value = 4070.3333
$ws.Cells.Item(132, 11).Value = 12210.9999
$ws.Cells.Item(132, 13).Value = -9680.999899999999


$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Cells.Item(2, 8).Value = 278.11765
$ws.Cells.Item(2, 9).Value = 143.09091
$ws.Cells.Item(2, 10).Value = 525.6667
$ws.Cells.Item(2, 11).Value = 858.54546
$ws.Cells.Item(2, 12).Value = 3154.0002
$ws.Cells.Item(2, 13).Value = -745.54546
$ws.Cells.Item(2, 14).Value = -3380.0002

# Row 16
$ws.Cells.Item(16, 8).Value = 326.66666
$ws.Cells.Item(16, 9).Value = 326.66666
$ws.Cells.Item(16, 11).Value = 979.9999799999999
$ws.Cells.Item(16, 13).Value = -806.9999799999999

# Row 104
$ws.Cells.Item(104, 8).Value = 1397.5
$ws.Cells.Item(104, 10).Value = 1395
$ws.Cells.Item(104, 12).Value = 4185
$ws.Cells.Item(104, 14).Value = -9427

# Row 108
$ws.Cells.Item(108, 8).Value = 484.83334
$ws.Cells.Item(108, 9).Value = 484.83334
$ws.Cells.Item(108, 11).Value = 1454.50002
$ws.Cells.Item(108, 13).Value = 1425.49998

# Row 109
$ws.Cells.Item(109, 8).Value = 985.4
$ws.Cells.Item(109, 9).Value = 985.4
$ws.Cells.Item(109, 10).Value = 0
$ws.Cells.Item(109, 11).Value = 2956.2
$ws.Cells.Item(109, 12).Value = 0
$ws.Cells.Item(109, 13).Value = -1916.2
$ws.Cells.Item(109, 14).ClearContents()

# Row 121
$ws.Cells.Item(121, 8).Value = 834.8570999999999
$ws.Cells.Item(121, 10).Value = 1161
$ws.Cells.Item(121, 12).Value = 3483
$ws.Cells.Item(121, 14).Value = -6103

# Row 128
$ws.Cells.Item(128, 8).Value = 301000.66
$ws.Cells.Item(128, 9).Value = 301000.66
$ws.Cells.Item(128, 11).Value = 903001.98
$ws.Cells.Item(128, 13).Value = -898021.98


$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Cells.Item(102, 8).Value = 11285.333
$ws.Cells.Item(102, 9).Value = 2300.818
$ws.Cells.Item(102, 11).Value = 2300.818
$ws.Cells.Item(102, 13).Value = -678.8180000000002

# Row 105
$ws.Cells.Item(105, 8).Value = 45202.855
$ws.Cells.Item(105, 9).Value = 18000
$ws.Cells.Item(105, 10).Value = 72405.71000000001
$ws.Cells.Item(105, 11).Value = 18000
$ws.Cells.Item(105, 12).Value = 72405.71000000001
$ws.Cells.Item(105, 13).Value = -14506
$ws.Cells.Item(105, 14).Value = -79393.71000000001

# Row 122
$ws.Cells.Item(122, 8).Value = 2551.1052
$ws.Cells.Item(122, 9).Value = 2309.5
$ws.Cells.Item(122, 10).Value = 2819.5557
$ws.Cells.Item(122, 11).Value = 6928.5
$ws.Cells.Item(122, 12).Value = 8458.667099999999
$ws.Cells.Item(122, 13).Value = -4478.5
$ws.Cells.Item(122, 14).Value = -13358.6671

# Row 132
$ws.Cells.Item(132, 8).Value = 2062
$ws.Cells.Item(132, 9).Value = 1876.7778
$ws.Cells.Item(132, 11).Value = 5630.3334
$ws.Cells.Item(132, 13).Value = -3100.3334


$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Cells.Item(40, 8).Value = 7342.5713
$ws.Cells.Item(40, 9).Value = 5350
$ws.Cells.Item(40, 10).Value = 9999.333000000001
$ws.Cells.Item(40, 11).Value = 5350
$ws.Cells.Item(40, 12).Value = 9999.333000000001
$ws.Cells.Item(40, 13).Value = -5214
$ws.Cells.Item(40, 14).Value = -10271.333

# Row 132
$ws.Cells.Item(132, 8).Value = 3000
$ws.Cells.Item(132, 9).Value = 1000
$ws.Cells.Item(132, 10).Value = 5000
$ws.Cells.Item(132, 11).Value = 3000
$ws.Cells.Item(132, 12).Value = 15000
$ws.Cells.Item(132, 13).Value = -470
$ws.Cells.Item(132, 14).Value = -20060

# Row 136
$ws.Cells.Item(136, 8).Value = 3797.6
$ws.Cells.Item(136, 9).Value = 3498
$ws.Cells.Item(136, 10).Value = 4247
$ws.Cells.Item(136, 11).Value = 10494
$ws.Cells.Item(136, 12).Value = 12741
$ws.Cells.Item(136, 13).Value = -7944
$ws.Cells.Item(136, 14).Value = -17841


$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Cells.Item(113, 8).Value = 394
$ws.Cells.Item(113, 9).Value = 394
$ws.Cells.Item(113, 11).Value = 1182
$ws.Cells.Item(113, 13).Value = 988

# Row 132
$ws.Cells.Item(132, 8).Value = 2377.25
$ws.Cells.Item(132, 9).Value = 1752
$ws.Cells.Item(132, 10).Value = 3002.5
$ws.Cells.Item(132, 11).Value = 5256
$ws.Cells.Item(132, 12).Value = 9007.5
$ws.Cells.Item(132, 13).Value = -2726
$ws.Cells.Item(132, 14).Value = -14067.5

# Row 136
$ws.Cells.Item(136, 8).Value = 6788.7856
$ws.Cells.Item(136, 9).Value = 5237.375
$ws.Cells.Item(136, 10).Value = 8857.333000000001
$ws.Cells.Item(136, 11).Value = 15712.125
$ws.Cells.Item(136, 12).Value = 26571.999
$ws.Cells.Item(136, 13).Value = -13162.125
$ws.Cells.Item(136, 14).Value = -31671.999

